$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# --- Copy existing formatting down onto the new cells before writing values ---
# Column B (date style, s=1) and Column C (time style, s=2) for rows 18-19
$ws.Range("B17").Copy()
$ws.Range("B18:B19").PasteSpecial($xlPasteFormats)
$ws.Range("C17").Copy()
$ws.Range("C18:C19").PasteSpecial($xlPasteFormats)

# Column D/E (time style, s=2) for row 17 (new D/E cells) and row 18 (new row)
$ws.Range("D16").Copy()
$ws.Range("D17:D18").PasteSpecial($xlPasteFormats)
$ws.Range("E16").Copy()
$ws.Range("E17:E18").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# --- Row 17: fill in the newly added D17 / E17 ---
$ws.Range("D17").Value = 0.59375
$ws.Range("E17").Formula = "=D17-C17"

# --- Row 18 (new row): Dataset entry ---
$ws.Range("B18").Value = 45999
$ws.Range("C18").Value = 0.625
$ws.Range("D18").Value = 0.64583333333333337
$ws.Range("E18").Formula = "=D18-C18"
$ws.Range("F18").Value = "Dataset"

# --- Row 19 (new row): Baseline Training entry (only B, C, F filled) ---
$ws.Range("B19").Value = 45999
$ws.Range("C19").Value = 0.64583333333333337
$ws.Range("F19").Value = "Baseline Training"

# --- Update selection / scroll position to match final state ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F19").Select()
